# Swap the species-record data between row 17 and row 18 on the "Artfynd"
# sheet. Columns P onward are identical between the two rows, so only the
# record-specific columns (A, B, D, E, F, G, H, M) need to be exchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Artfynd")

$cols = @("A", "B", "D", "E", "F", "G", "H", "M")

foreach ($col in $cols) {
    $cell17 = $ws.Range($col + "17")
    $cell18 = $ws.Range($col + "18")

    $val17 = $cell17.Value2
    $val18 = $cell18.Value2

    $cell17.Value = $val18
    $cell18.Value = $val17
}
